$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows right after the existing row 113 (i.e. at 114-115),
# pushing the previous rows 114-193 down to 116-195.
$ws.Rows("114:115").Insert()

# --- New row 114 ---
$ws.Range("A114").Value = 1
$ws.Range("B114").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C114").Value = 'Arica y Parinacota'
$ws.Range("D114").Value = 45126
$ws.Range("E114").Value = 15
$ws.Range("F114").Value = 'Fruta'
$ws.Range("G114").Value = 100108
$ws.Range("H114").Value = 'Tropicales y subtropicales'
$ws.Range("I114").Value = 100108003
$ws.Range("J114").Value = 'Maracuyá'
$ws.Range("K114").Value = 'Sin especificar'
$ws.Range("L114").Value = 'Primera'
$ws.Range("M114").Value = 130
$ws.Range("N114").Value = 24000
$ws.Range("O114").Value = 25000
$ws.Range("P114").Value = 24462
$ws.Range("Q114").Value = '$/caja 20 kilos'
$ws.Range("R114").Value = 'Región de Arica y Parinacota'
$ws.Range("S114").Value = 1223
$ws.Range("T114").Value = 20

# --- New row 115 ---
$ws.Range("A115").Value = 1
$ws.Range("B115").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C115").Value = 'Arica y Parinacota'
$ws.Range("D115").Value = 45126
$ws.Range("E115").Value = 15
$ws.Range("F115").Value = 'Fruta'
$ws.Range("G115").Value = 100108
$ws.Range("H115").Value = 'Tropicales y subtropicales'
$ws.Range("I115").Value = 100108003
$ws.Range("J115").Value = 'Maracuyá'
$ws.Range("K115").Value = 'Sin especificar'
$ws.Range("L115").Value = 'Segunda'
$ws.Range("M115").Value = 180
$ws.Range("N115").Value = 20000
$ws.Range("O115").Value = 21000
$ws.Range("P115").Value = 20556
$ws.Range("Q115").Value = '$/caja 20 kilos'
$ws.Range("R115").Value = 'Región de Arica y Parinacota'
$ws.Range("S115").Value = 1028
$ws.Range("T115").Value = 20
